$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DeliveryDate cell (B2) used to hold a numeric date value formatted as
# DD/MM/YYYY. Replace it with the literal text "31-12-2018" (a string,
# not a real date), keeping its existing style/number format.
$ws.Range("B2").Value = "31-12-2018"

# Update the active selection to B2, matching the saved view state.
$ws.Range("B2").Select()
